# feat: add 2022-Q1 data
#
# The workbook previously ended with a single summary sheet "总计".
# This edit:
#   1. Turns that sheet into the new per-quarter holdings sheet "2022-Q1"
#      (same shape/format as the other quarter sheets).
#   2. Appends a brand-new "总计" sheet at the end containing the original
#      summary rows plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet into "2022-Q1"
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# A style-2 cell (bold, centered, bordered) already exists elsewhere in the
# workbook (e.g. the other quarter sheets) - copy its format instead of
# re-building it from scratch so we don't create a duplicate style record.
$styleSrc = $wb.Worksheets.Item("2021-Q4").Range("A2")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q1.Cells.Item(1, $col).Value2 = $headers[$col - 2]
}

$fundRows = @(
    @("010379", "广发均衡优选混合A", "55.11", "64.69", "4.95", "2.7279", 7),
    @("001532", "华安文体健康主题灵活配置混合", "87.24", "92.86", "2.11", "1.8408", 9),
    @("009887", "广发稳健优选六个月持有期混合A", "30.91", "64.69", "5.06", "1.5640", 5),
    @("011194", "广发睿铭两年持有期混合型证券投资基金A", "20.55", "74.56", "4.89", "1.0049", 8),
    @("002350", "华安安华灵活配置混合", "42.47", "93.61", "2.15", "0.9131", 8),
    @("009888", "广发稳健优选六个月持有期混合C", "11.86", "64.69", "5.06", "0.6001", 5),
    @("011195", "广发睿铭两年持有期混合型证券投资基金C", "7.63", "74.56", "4.89", "0.3731", 8),
    @("011134", "广发价值优选混合A", "6.35", "93.95", "4.92", "0.3124", 10),
    @("010380", "广发均衡优选混合C", "3.49", "64.69", "4.95", "0.1728", 7),
    @("008531", "惠升惠民混合A", "3.71", "71.70", "2.14", "0.0794", 10),
    @("011135", "广发价值优选混合C", "1.48", "93.95", "4.92", "0.0728", 10),
    @("080005", "长盛量化红利混合", "2.66", "69.88", "2.52", "0.0670", 8),
    @("008532", "惠升惠民混合C", "1.20", "71.70", "2.14", "0.0257", 10)
)

$row = 2
foreach ($r in $fundRows) {
    $q1.Cells.Item($row, 1).Value2 = $row - 2
    $q1.Cells.Item($row, 2).Value2 = $r[0]
    $q1.Cells.Item($row, 3).Value2 = $r[1]
    $q1.Cells.Item($row, 4).Value2 = $r[2]
    $q1.Cells.Item($row, 5).Value2 = $r[3]
    $q1.Cells.Item($row, 6).Value2 = $r[4]
    $q1.Cells.Item($row, 7).Value2 = $r[5]
    $q1.Cells.Item($row, 8).Value2 = $r[6]
    $row = $row + 1
}

# Re-apply the numeric-looking text columns (B, D, E, F, G) as text so the
# leading zeros / trailing zeros of fund codes and ratios are preserved,
# matching how the other quarter sheets store this data.
$q1.Range("B2:B14").NumberFormat = "@"
for ($r = 2; $r -le 14; $r++) {
    $q1.Cells.Item($r, 2).Value2 = $fundRows[$r - 2][0]
}
$q1.Range("D2:G14").NumberFormat = "@"
for ($r = 2; $r -le 14; $r++) {
    $rowData = $fundRows[$r - 2]
    $q1.Cells.Item($r, 4).Value2 = $rowData[2]
    $q1.Cells.Item($r, 5).Value2 = $rowData[3]
    $q1.Cells.Item($r, 6).Value2 = $rowData[4]
    $q1.Cells.Item($r, 7).Value2 = $rowData[5]
}

# Stamp the shared bold/centered/bordered style onto the header row and the
# index column, exactly like the other quarter sheets.
$styleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet with the refreshed summary table
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

$total.Cells.Item(1, 2).Value2 = "日期"
$total.Cells.Item(1, 3).Value2 = "持有数量(只)"
$total.Cells.Item(1, 4).Value2 = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 13, 9.75),
    @("2021-Q4", 20, 5.02),
    @("2021-Q3", 7, 1.71),
    @("2021-Q1", 2, 0),
    @("2020-Q4", 3, 0.23)
)

$row = 2
foreach ($r in $summaryRows) {
    $total.Cells.Item($row, 1).Value2 = $row - 2
    $total.Cells.Item($row, 2).Value2 = $r[0]
    $total.Cells.Item($row, 3).Value2 = $r[1]
    $total.Cells.Item($row, 4).Value2 = $r[2]
    $row = $row + 1
}

$styleSrc.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the original active sheet/selection (the workbook was focused on
# "2020-Q4" before this edit and none of the sheet-addition/renaming above
# should change that).
$wb.Worksheets.Item("2020-Q4").Activate()
